# Replace the monthly-resetting index in column A with a continuously
# incrementing index (0-based) running across the whole dataset.
# Rows 2..23 (January) already hold the correct values (0..21) and are
# left untouched; rows 24..248 (February onward) are rewritten so the
# counter keeps climbing instead of resetting to 0 at the start of each
# month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 24
$endRow = 248
$counter = $startRow - 2   # continuous index value for $startRow (22)

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $counter
    $counter++
}
